$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the Correspond Handoff / Handback datetimes for the
# "41246388-cc5f-4bea-ad44-58465805de5c" row (row 3)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-23 12:48:38"
$wsZh.Range("H3").Value = "2016-03-23 12:49:06"

# de-de sheet: update the Correspond Handoff / Handback datetimes for the
# "41246388-cc5f-4bea-ad44-58465805de5c" row (row 3)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-23 12:48:43"
$wsDe.Range("H3").Value = "2016-03-23 12:49:12"
